# Update the SDG 12.4.2 metadata sheet with the newly uploaded indicator
# passport data (re-upload of meta/12-4-2.xlsx).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: Indicator wording was refined (added a)/b) split) ---
$ws.Range("B4").Value = "12.4.2 a) Образование опасных отходов на душу населения и b) доля обрабатываемых опасных отходов в разбивке по видам обработки "

# --- Section 2: Organization / contact info was updated ---
$ws.Range("B6").Value = "Национальный статистический комитет КР (Управление цифрового развития и статистики устойчивого развития)"
$ws.Range("B7").Value = "Мамбеталиев Т.А."
$ws.Range("B9").Value = "(0312) 62 56 07"
$ws.Range("B10").Value = "www.stat.gov.kg"

# --- Selection left on the contact-person cell, as in the saved file ---
$ws.Range("B7").Select()
